# Update the click/impression data on the ClickThroughRateSheet.
# Only the raw values in columns B (clicks) and C (impressions) for rows 3-13
# change; column D holds the shared CTR formula and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClickThroughRateSheet")

$newValues = @{
    3  = @{ B = 590; C = 614036 }
    4  = @{ B = 795; C = 200127 }
    5  = @{ B = 858; C = 341388 }
    6  = @{ B = 849; C = 408012 }
    7  = @{ B = 892; C = 202438 }
    8  = @{ B = 847; C = 811081 }
    9  = @{ B = 391; C = 377913 }
    10 = @{ B = 854; C = 157944 }
    11 = @{ B = 643; C = 877947 }
    12 = @{ B = 378; C = 388269 }
    13 = @{ B = 643; C = 747905 }
}

foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value = $newValues[$row].B
    $ws.Range("C$row").Value = $newValues[$row].C
}
